# Apply the "food patchwork" (previously only applied to the Electricity/heat
# column) to the remaining non-food sectors in rows 2 (Agriculture), 3 (Animals)
# and 9 (Consumable): flip the relevant cells from 1 (white) to 0 (red fill),
# matching the pattern already used for the Electricity/heat column (style/
# fill index carried by the red "0" cells elsewhere in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Agriculture
$ws.Range("E2:H2").Value = 0
$ws.Range("E2:H2").Interior.Color = 255
$ws.Range("J2:S2").Value = 0
$ws.Range("J2:S2").Interior.Color = 255
$ws.Range("U2:X2").Value = 0
$ws.Range("U2:X2").Interior.Color = 255

# Row 3 - Animals
$ws.Range("E3:H3").Value = 0
$ws.Range("E3:H3").Interior.Color = 255
$ws.Range("J3:K3").Value = 0
$ws.Range("J3:K3").Interior.Color = 255
$ws.Range("M3:S3").Value = 0
$ws.Range("M3:S3").Interior.Color = 255
$ws.Range("U3:X3").Value = 0
$ws.Range("U3:X3").Interior.Color = 255

# Row 9 - Consumable
$ws.Range("B9:H9").Value = 0
$ws.Range("B9:H9").Interior.Color = 255
$ws.Range("J9:S9").Value = 0
$ws.Range("J9:S9").Interior.Color = 255
$ws.Range("U9:X9").Value = 0
$ws.Range("U9:X9").Interior.Color = 255

# Move the active selection, as recorded in the saved workbook view.
$ws.Range("M26").Select() | Out-Null
